# Auto-generated edit script: updates Leve profit-calculation cells
# across all 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to
# reflect refreshed market-board price data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 82.8
$ws.Range("I15").Value = 82.8
$ws.Range("K15").Value = 248.4
$ws.Range("M15").Value = -79.39999999999998
$ws.Range("H132").Value = 25927.795
$ws.Range("I132").Value = 3875.8647
$ws.Range("J132").Value = 433888.5
$ws.Range("K132").Value = 11627.5941
$ws.Range("L132").Value = 1301665.5
$ws.Range("M132").Value = -9097.5941
$ws.Range("N132").Value = -1306725.5
$ws.Range("H135").Value = 27778780
$ws.Range("I135").Value = 1060.4117
$ws.Range("K135").Value = 9543.705300000001
$ws.Range("M135").Value = -7008.705300000001
$ws.Range("H138").Value = 2775.25
$ws.Range("I138").Value = 2330.423
$ws.Range("J138").Value = 3026.6738
$ws.Range("K138").Value = 6991.268999999999
$ws.Range("L138").Value = 9080.0214
$ws.Range("M138").Value = -1851.268999999999
$ws.Range("N138").Value = -19360.0214

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30985.705
$ws.Range("I32").Value = 31168.8
$ws.Range("K32").Value = 31168.8
$ws.Range("M32").Value = -30881.8
$ws.Range("H61").Value = 3526
$ws.Range("I61").Value = 3529.2917
$ws.Range("J61").Value = 3499.6667
$ws.Range("K61").Value = 3529.2917
$ws.Range("L61").Value = 3499.6667
$ws.Range("M61").Value = -3317.2917
$ws.Range("N61").Value = -3923.6667
$ws.Range("H109").Value = 41374.332
$ws.Range("J109").Value = 41374.332
$ws.Range("L109").Value = 41374.332
$ws.Range("N109").Value = -44148.332
$ws.Range("H122").Value = 1096.8334
$ws.Range("I122").Value = 1016.3
$ws.Range("J122").Value = 1499.5
$ws.Range("K122").Value = 3048.9
$ws.Range("L122").Value = 4498.5
$ws.Range("M122").Value = -598.8999999999996
$ws.Range("N122").Value = -9398.5
$ws.Range("H132").Value = 17244024
$ws.Range("I132").Value = 21741370
$ws.Range("J132").Value = 4199
$ws.Range("K132").Value = 65224110
$ws.Range("L132").Value = 12597
$ws.Range("M132").Value = -65221580
$ws.Range("N132").Value = -17657
$ws.Range("H136").Value = 3526
$ws.Range("I136").Value = 3529.2917
$ws.Range("J136").Value = 3499.6667
$ws.Range("K136").Value = 10587.8751
$ws.Range("L136").Value = 10499.0001
$ws.Range("M136").Value = -8037.875100000001
$ws.Range("N136").Value = -15599.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2702.8667
$ws.Range("I134").Value = 1794.7
$ws.Range("J134").Value = 4519.2
$ws.Range("K134").Value = 5384.1
$ws.Range("L134").Value = 13557.6
$ws.Range("M134").Value = -2849.1
$ws.Range("N134").Value = -18627.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10400.706
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 10400.706
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 10400.706
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = -10990.706
$ws.Range("H34").Value = 10400.706
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 10400.706
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 10400.706
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -10804.706
$ws.Range("H132").Value = 66563.73
$ws.Range("I132").Value = 2191.6924
$ws.Range("K132").Value = 6575.0772
$ws.Range("M132").Value = -4045.0772
$ws.Range("H134").Value = 3925.476
$ws.Range("I134").Value = 1713.6666
$ws.Range("J134").Value = 5584.3335
$ws.Range("K134").Value = 5140.9998
$ws.Range("L134").Value = 16753.0005
$ws.Range("M134").Value = -2605.9998
$ws.Range("N134").Value = -21823.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 13169.8125
$ws.Range("I107").Value = 10651.3
$ws.Range("J107").Value = 17367.334
$ws.Range("K107").Value = 31953.9
$ws.Range("L107").Value = 52102.00199999999
$ws.Range("M107").Value = -30033.9
$ws.Range("N107").Value = -55942.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 13139.637
$ws.Range("I46").Value = 28978
$ws.Range("J46").Value = 9620
$ws.Range("K46").Value = 28978
$ws.Range("L46").Value = 9620
$ws.Range("M46").Value = -28822
$ws.Range("N46").Value = -9932
$ws.Range("H57").Value = 34526.25
$ws.Range("J57").Value = 36744.285
$ws.Range("L57").Value = 36744.285
$ws.Range("N57").Value = -38384.285
$ws.Range("H80").Value = 195861.53
$ws.Range("I80").Value = 316337.5
$ws.Range("J80").Value = 3100
$ws.Range("K80").Value = 316337.5
$ws.Range("L80").Value = 3100
$ws.Range("M80").Value = -315339.5
$ws.Range("N80").Value = -5096
$ws.Range("H83").Value = 195861.53
$ws.Range("I83").Value = 316337.5
$ws.Range("J83").Value = 3100
$ws.Range("K83").Value = 1581687.5
$ws.Range("L83").Value = 15500
$ws.Range("M83").Value = -1576695.5
$ws.Range("N83").Value = -25484
$ws.Range("H102").Value = 1051.4706
$ws.Range("I102").Value = 731.8182
$ws.Range("J102").Value = 1637.5
$ws.Range("K102").Value = 731.8182
$ws.Range("L102").Value = 1637.5
$ws.Range("M102").Value = 890.1818
$ws.Range("N102").Value = -4881.5
$ws.Range("H132").Value = 3357.8333
$ws.Range("I132").Value = 2779.077
$ws.Range("J132").Value = 4041.818
$ws.Range("K132").Value = 8337.231
$ws.Range("L132").Value = 12125.454
$ws.Range("M132").Value = -5807.231
$ws.Range("N132").Value = -17185.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 2822
$ws.Range("I9").Value = 225.14285
$ws.Range("K9").Value = 225.14285
$ws.Range("M9").Value = -1.14285000000001
$ws.Range("H132").Value = 3997.5
$ws.Range("I132").Value = 3771.52
$ws.Range("J132").Value = 4374.1333
$ws.Range("K132").Value = 11314.56
$ws.Range("L132").Value = 13122.3999
$ws.Range("M132").Value = -8784.559999999999
$ws.Range("N132").Value = -18182.3999
$ws.Range("H136").Value = 2683.6072
$ws.Range("I136").Value = 2254.05
$ws.Range("J136").Value = 3757.5
$ws.Range("K136").Value = 6762.150000000001
$ws.Range("L136").Value = 11272.5
$ws.Range("M136").Value = -4212.150000000001
$ws.Range("N136").Value = -16372.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 19577480
$ws.Range("I122").Value = 31093144
$ws.Range("J122").Value = 851.9
$ws.Range("K122").Value = 93279432
$ws.Range("L122").Value = 2555.7
$ws.Range("M122").Value = -93276982
$ws.Range("N122").Value = -7455.7
$ws.Range("H132").Value = 1709.742
$ws.Range("I132").Value = 1456.6522
$ws.Range("J132").Value = 2437.375
$ws.Range("K132").Value = 4369.9566
$ws.Range("L132").Value = 7312.125
$ws.Range("M132").Value = -1839.9566
$ws.Range("N132").Value = -12372.125
